# cdm324_v2 BOM update — "first stage changed, POR for segment controller"
#
# - C4 (row 3): value/manufacturer/part-number/LCSC-part/link updated from the
#   old 200p Samsung cap to a 27pF YAGEO cap.
# - FB1 (row 6): Package/Footprint filled in as 0805.
# - R3, R5, R6 (row 11) split: that row now only covers R5/R6 (qty 3 -> 2).
# - New row 18 added for R3 on its own (620k, YAGEO RC0603FR-07620KL, 0603).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Item 2, designator C4): part swap ---
$ws.Range("D3").Value = "27pF"
$ws.Range("E3").Value = "YAGEO"
$ws.Range("F3").Value = "CC0603JRNPO9BN270"
$ws.Range("G3").Value = "C107045"

# --- Row 6 (Item 5, designator FB1): fill in Package/Footprint ---
# Leading-zero text value - quote-prefix so it is stored as text (not 805).
$ws.Range("H6").Value = "'0805"

# --- Row 11 (Item 10): designator list shrinks to R5, R6 only; qty 3 -> 2 ---
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "R5, R6"

# --- New row 18 (Item 17): R3 now its own BOM line ---
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "R3"
$ws.Range("D18").Value = "620k"
$ws.Range("E18").Value = "YAGEO"
$ws.Range("F18").Value = "RC0603FR-07620KL"
$ws.Range("G18").Value = "C68737"
$ws.Range("H18").Value = "'0603"

# --- Hyperlinks: the I3 (C4 Digikey link) target must point at the new part.
# This runtime's Range.Hyperlinks.Delete() removes every hyperlink on the
# sheet (not just the targeted range), so capture all existing targets first,
# wipe them, then re-add all of them with I3 updated to the new URL.
$linkOrder = @("I9","I2","I4","I5","I6","I7","I8","I10","I11","I12","I13","I14","I15","I16","I17","I3")
$linkTargets = @{
  "I2"  = "https://www.digikey.com/product-detail/en/samsung-electro-mechanics/CL21A226MOCLRNC/1276-6780-1-ND/5961639"
  "I3"  = "https://www.digikey.com/en/products/detail/yageo/CC0603JRNPO9BN270/302797"
  "I4"  = "https://www.digikey.com/en/products/detail/tdk-corporation/C1005X5R1E474K050BB/2792240"
  "I5"  = "https://www.digikey.com/en/products/detail/murata-electronics/GRM188R71A225KE15D/1033256"
  "I6"  = "http://www.digikey.com/product-detail/en/HZ0805D102R-10/240-2397-1-ND/806757"
  "I7"  = "http://www.digikey.com/product-detail/en/sullins-connector-solutions/PPPC031LFBN-RC/S7036-ND/810175"
  "I8"  = "https://www.digikey.com/en/products/detail/sullins-connector-solutions/PPPC042LFBN-RC/810244"
  "I9"  = "https://www.digikey.com/en/products/detail/yageo/RC0603FR-074K7L/727212"
  "I10" = "https://www.digikey.com/en/products/detail/yageo/RC0402JR-0775RL/726505"
  "I11" = "https://www.digikey.com/product-detail/en/panasonic-electronic-components/ERJ-3EKF1303V/P130KHCT-ND/198154"
  "I12" = "https://www.digikey.ch/en/products/detail/texas-instruments/LMP7731MF-NOPB/1836050"
  "I13" = "https://www.digikey.com/en/products/detail/analog-devices-inc-maxim-integrated/MAX9814ETD-T/1703775"
  "I14" = "https://www.digikey.com/product-detail/en/on-semiconductor/NCP163ASN330T1G/NCP163ASN330T1GOSCT-ND/10233513"
  "I15" = "https://www.digikey.com/product-detail/en/stmicroelectronics/STM32F301K8U6/497-17411-ND/5051319"
  "I16" = "https://www.digikey.ch/en/products/detail/nxp-usa-inc/PCA8561AHN-AY/5170024"
  "I17" = "https://www.aliexpress.com/item/1005003060952508.html"
}

[void]$ws.Range("I2").Hyperlinks.Delete()
foreach ($addr in $linkOrder) {
  [void]$ws.Hyperlinks.Add($ws.Range($addr), $linkTargets[$addr])
}

# The displayed cell text mirrors the URL in this BOM - update it to match
# the new Digikey link now that the hyperlink target itself changed.
$ws.Range("I3").Value = "https://www.digikey.com/en/products/detail/yageo/CC0603JRNPO9BN270/302797"

# Restore the active-cell selection to match the saved view.
[void]$ws.Range("D19").Select()
